$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.690.56'
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").Value = '3.262.62'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.63'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.36'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.129'
$ws.Range("E9").Value = '  -4.14%  '
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("E11").Value = '  -3.54%  '
$ws.Range("D12").Value = '3.824.30'
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.33'
$ws.Range("E14").Value = '  -5.82%  '
$ws.Range("D15").Value = '67.748.68'
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000167'
$ws.Range("E16").Value = '  -2.72%  '
$ws.Range("D17").Value = '3.252.17'
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.70'
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.41'
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '394.68'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.54'
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.76'
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("E25").Value = '  -4.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.187'
$ws.Range("E26").Value = '  -0.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.46'
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("E29").Value = '  -2.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.56'
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.44'
$ws.Range("E31").Value = '  -5.70%  '
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  -5.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.15'
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("E36").Value = '  -5.54%  '
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.75'
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("E39").Value = '  -3.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.48'
$ws.Range("E40").Value = '  -2.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.27'
$ws.Range("E41").Value = '  -4.79%  '
$ws.Range("D42").Value = '2.655.42'
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.77'
$ws.Range("E43").Value = '  -1.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0679'
$ws.Range("E44").Value = '  -2.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  -6.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.58'
$ws.Range("E46").Value = '  -2.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '334.59'
$ws.Range("E47").Value = '  -2.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0273'
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.30'
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.100'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("E51").Value = '  -2.26%  '
